# edit.ps1 - applies the changes described by the diff to Vinh H Nguyen.docx
# Runs against the Word COM object model (iron_native headless runtime).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0. Remove the stray "_GoBack" bookmark from its old location (near the
#    "Master of Social Research ... Macquarie Australia" / Australia
#    address block). It will be re-added at its new location later
#    (right after "BLUEPRISM/PEGA" in the skills bullet inside the table).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1. AIG bullet: rewrite the "Led the testing..." sentence.
# ---------------------------------------------------------------------
$old1 = "Led the testing and collaboration with external auditor on the SOX auditing of IT General Control and Application Controls at AIG Finance, Investments and Corporate"
$new1 = "Led the SOX auditing (methodology and audit findings fully disclosed to external auditors) of ITGC general control, IT dependent and Application Controls for Finance and Investments "
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------
# 2/3. Both "BLUEPRISM RPA, " occurrences -> "BLUEPRISM/PEGA RPA, ":
#      - one inside the table (skills bullet list)
#      - one outside the table (further down, followed by "APPIAN ...")
#    A single ReplaceAll catches both (Find text can span run
#    boundaries), then the "_GoBack" bookmark is re-inserted right
#    after the newly added "/PEGA" at the FIRST (table) occurrence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("BLUEPRISM RPA, ", $false, $false, $false, $false, $false, $true, 1, $false, "BLUEPRISM/PEGA RPA, ", 2)

$bpFind = $d.Content
$bpFind.Find.Execute("BLUEPRISM/PEGA")
$bpEnd = $bpFind.End
$bpBookmarkRange = $d.Range($bpEnd, $bpEnd)
$d.Bookmarks.Add("_GoBack", $bpBookmarkRange)

# ---------------------------------------------------------------------
# 4. Education bullet: " Master of Social Research, Macquarie Australia"
#    -> " Master of Social Research (Macquarie Australia)".
# ---------------------------------------------------------------------
$old4 = " Master of Social Research, Macquarie Australia"
$new4 = " Master of Social Research (Macquarie Australia)"
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2)

# ---------------------------------------------------------------------
# 5. Append "S" to the two "EXPERIENCE" headings so they read
#    "...EXPERIENCES".
#    a) "SELECTED PROFESSIONAL EXPERIENCE" (first occurrence)
#    b) "GLOBAL EXPERIENCE" (second occurrence)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("SELECTED PROFESSIONAL EXPERIENCE", $false, $false, $false, $false, $false, $true, 1, $false, "SELECTED PROFESSIONAL EXPERIENCES", 2)
$d.Content.Find.Execute("GLOBAL EXPERIENCE", $false, $false, $false, $false, $false, $true, 1, $false, "GLOBAL EXPERIENCES", 2)

# ---------------------------------------------------------------------
# 6. Rewrite the "Developed capabilities for 100% internal reliance..."
#    sentence.
# ---------------------------------------------------------------------
$old6 = "for 100% internal reliance in the SOX/MAR documentation and testing of 20+ controls and 50 test plans across 150 applications in Finance, Actuarial, Risk, Treasury " + [char]8211 + " with close relationship to "
$new6 = "for 100% insourcing SOX/MAR auditing including the documentation and testing of controls across 150 applications in Finance, Treasury, Actuarial and Risk.  Audit methodology and findings  are fully shared with "
$d.Content.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
